# Generate Report for Handback
# Adds a new handback row (e48a681a-322c-41fe-9339-23f16b108803.md) to the
# Overview / zh-cn / de-de sheets, matching the target OOXML diff.

$wb = $excel.ActiveWorkbook

$HYPER_COLOR = 15570276  # RGB(100,149,237) == #6495ED stored as 0xBBGGRR for Excel Font.Color

function Set-HyperlinkLook($rng) {
    # Re-create the workbook's "HyperLink" look (blue + underline) on a cell
    $rng.Font.Color = $HYPER_COLOR
    $rng.Font.Underline = 2
}

function Set-DateLook($rng) {
    $rng.NumberFormat = "yyyy-mm-dd HH:mm:ss"
}

# ---------------------------------------------------------------------------
# Sheet "Overview"  (row 4)
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A4").Value = "e48a681a-322c-41fe-9339-23f16b108803.md"

$wsOverview.Range("B4").Value = "e2e\e48a681a-322c-41fe-9339-23f16b108803.md"
$wsOverview.Hyperlinks.Add(
    $wsOverview.Range("B4"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/15160d73d0e7c26470bcbf9acc29108c8c7004f5/e2e/e48a681a-322c-41fe-9339-23f16b108803.md",
    [Type]::Missing,
    [Type]::Missing,
    "e2e\e48a681a-322c-41fe-9339-23f16b108803.md"
) | Out-Null
Set-HyperlinkLook $wsOverview.Range("B4")

$wsOverview.Range("C4").Value = ".md"
$wsOverview.Range("E4").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F4").Value = "Handed back: in sync with en-US"

$wsOverview.Range("G4").Value = "2016-08-26 20:44:13"
Set-DateLook $wsOverview.Range("G4")

$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.Resize($wsOverview.Range("A1:G4"))

# ---------------------------------------------------------------------------
# Sheet "zh-cn"  (row 4)
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("A4").Value = "e48a681a-322c-41fe-9339-23f16b108803.md"
$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("A4"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/583b319067623eadc2e9547302921c723e0ff1c2/e2e/e48a681a-322c-41fe-9339-23f16b108803.md",
    [Type]::Missing,
    [Type]::Missing,
    "e48a681a-322c-41fe-9339-23f16b108803.md"
) | Out-Null
Set-HyperlinkLook $wsZhCn.Range("A4")

$wsZhCn.Range("B4").Value = ".md"
$wsZhCn.Range("C4").Value = "Handed back: in sync with en-US"
$wsZhCn.Range("D4").Value = "e2e"
$wsZhCn.Range("E4").Value = "ht"
$wsZhCn.Range("F4").Value = "True"
$wsZhCn.Range("G4").Value = "e48a681a-322c-41fe-9339-23f16b108803.201728c8f308229cca4c72b5a0f83b1b3dae8dcf.zh-cn.xlf"

$wsZhCn.Range("H4").Value = "2016-08-26 20:44:06"
Set-DateLook $wsZhCn.Range("H4")

$wsZhCn.Range("I4").Value = "e48a681a-322c-41fe-9339-23f16b108803.md"
$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("I4"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/583b319067623eadc2e9547302921c723e0ff1c2/e2e/e48a681a-322c-41fe-9339-23f16b108803.md",
    [Type]::Missing,
    [Type]::Missing,
    "e48a681a-322c-41fe-9339-23f16b108803.md"
) | Out-Null
Set-HyperlinkLook $wsZhCn.Range("I4")

$wsZhCn.Range("J4").Value = "e48a681a-322c-41fe-9339-23f16b108803.201728c8f308229cca4c72b5a0f83b1b3dae8dcf.zh-cn.xlf"

$wsZhCn.Range("K4").Value = "2016-08-26 20:44:28"
Set-DateLook $wsZhCn.Range("K4")

$wsZhCn.Range("L4").Value = ""
$wsZhCn.Range("M4").Value = "True"
$wsZhCn.Range("N4").Value = ""
$wsZhCn.Range("O4").Value = "False"
$wsZhCn.Range("P4").Value = ""

$loZhCn = $wsZhCn.ListObjects.Item(1)
$loZhCn.Resize($wsZhCn.Range("A1:P4"))

# ---------------------------------------------------------------------------
# Sheet "de-de"  (row 4)
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("A4").Value = "e48a681a-322c-41fe-9339-23f16b108803.md"
$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("A4"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/42e4c510de011f030ee169fe5a6f75a5e930a255/e2e/e48a681a-322c-41fe-9339-23f16b108803.md",
    [Type]::Missing,
    [Type]::Missing,
    "e48a681a-322c-41fe-9339-23f16b108803.md"
) | Out-Null
Set-HyperlinkLook $wsDeDe.Range("A4")

$wsDeDe.Range("B4").Value = ".md"
$wsDeDe.Range("C4").Value = "Handed back: in sync with en-US"
$wsDeDe.Range("D4").Value = "e2e"
$wsDeDe.Range("E4").Value = "ht"
$wsDeDe.Range("F4").Value = "True"
$wsDeDe.Range("G4").Value = "e48a681a-322c-41fe-9339-23f16b108803.201728c8f308229cca4c72b5a0f83b1b3dae8dcf.de-de.xlf"

$wsDeDe.Range("H4").Value = "2016-08-26 20:44:13"
Set-DateLook $wsDeDe.Range("H4")

$wsDeDe.Range("I4").Value = "e48a681a-322c-41fe-9339-23f16b108803.md"
$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("I4"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/42e4c510de011f030ee169fe5a6f75a5e930a255/e2e/e48a681a-322c-41fe-9339-23f16b108803.md",
    [Type]::Missing,
    [Type]::Missing,
    "e48a681a-322c-41fe-9339-23f16b108803.md"
) | Out-Null
Set-HyperlinkLook $wsDeDe.Range("I4")

$wsDeDe.Range("J4").Value = "e48a681a-322c-41fe-9339-23f16b108803.201728c8f308229cca4c72b5a0f83b1b3dae8dcf.de-de.xlf"

$wsDeDe.Range("K4").Value = "2016-08-26 20:44:35"
Set-DateLook $wsDeDe.Range("K4")

$wsDeDe.Range("L4").Value = ""
$wsDeDe.Range("M4").Value = "True"
$wsDeDe.Range("N4").Value = ""
$wsDeDe.Range("O4").Value = "False"
$wsDeDe.Range("P4").Value = ""

$loDeDe = $wsDeDe.ListObjects.Item(1)
$loDeDe.Resize($wsDeDe.Range("A1:P4"))

Write-Host "Handback report row added for e48a681a-322c-41fe-9339-23f16b108803.md"
